$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on price cells whose new values look numeric,
# so Excel keeps them as literal text (preserving formatting like trailing zeros)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values from the crypto data refresh
$ws.Range("D2").Value = "46.622.49"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.589.06"
$ws.Range("E3").Value = "  +10.63%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "305.21"
$ws.Range("D6").Value = "102.62"
$ws.Range("E6").Value = "  +4.05%  "
$ws.Range("D7").Value = "0.600"
$ws.Range("E7").Value = "  +5.90%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.568"
$ws.Range("E9").Value = "  +11.48%  "
$ws.Range("D10").Value = "38.61"
$ws.Range("E10").Value = "  +11.65%  "
$ws.Range("E11").Value = "  +5.09%  "
$ws.Range("E12").Value = "  +12.62%  "
$ws.Range("D13").Value = "2.984.67"
$ws.Range("E13").Value = "  +10.74%  "
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("D15").Value = "2.572.72"
$ws.Range("E15").Value = "  +10.03%  "
$ws.Range("E16").Value = "  +11.91%  "
$ws.Range("D17").Value = "15.03"
$ws.Range("E17").Value = "  +9.91%  "
$ws.Range("D18").Value = "47.587.01"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").Value = "13.92"
$ws.Range("E19").Value = "  +10.29%  "
$ws.Range("D20").Value = "0.0₂01000"
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("D21").Value = "6.58"
$ws.Range("E21").Value = "  +9.67%  "
$ws.Range("D22").Value = "69.94"
$ws.Range("E22").Value = "  +4.81%  "
$ws.Range("D23").Value = "254.58"
$ws.Range("E23").Value = "  +4.08%  "
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").Value = "  +6.14%  "
$ws.Range("E25").Value = "  +11.35%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "26.18"
$ws.Range("E27").Value = "  +24.32%  "
$ws.Range("D28").Value = "41.12"
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("D29").Value = "10.40"
$ws.Range("E29").Value = "  +7.95%  "
$ws.Range("D30").Value = "2.28"
$ws.Range("E30").Value = "  +3.30%  "
$ws.Range("E31").Value = "  +4.82%  "
$ws.Range("D32").Value = "2.94"
$ws.Range("E32").Value = "  +4.63%  "
$ws.Range("D33").Value = "5.96"
$ws.Range("E33").Value = "  +9.00%  "
$ws.Range("D34").Value = "0.0839"
$ws.Range("E34").Value = "  +8.77%  "
$ws.Range("E35").Value = "  +21.58%  "
$ws.Range("D36").Value = "148.29"
$ws.Range("E36").Value = "  +2.86%  "
$ws.Range("D37").Value = "0.120"
$ws.Range("E37").Value = "  +8.68%  "
$ws.Range("E38").Value = "  +3.87%  "
$ws.Range("D39").Value = "16.28"
$ws.Range("E39").Value = "  +8.98%  "
$ws.Range("D40").Value = "4.23"
$ws.Range("E40").Value = "  +9.65%  "
$ws.Range("D41").Value = "0.0327"
$ws.Range("E41").Value = "  +9.23%  "
$ws.Range("D42").Value = "3.58"
$ws.Range("E42").Value = "  +11.72%  "
$ws.Range("D43").Value = "2.023.35"
$ws.Range("E43").Value = "  +10.12%  "
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "92.63"
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("D46").Value = "17.67"
$ws.Range("E46").Value = "  +35.94%  "
$ws.Range("D47").Value = "1.85"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.199"
$ws.Range("E48").Value = "  +7.97%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "106.89"
$ws.Range("E49").Value = "  +10.72%  "
$ws.Range("D50").Value = "2.843.74"
$ws.Range("E50").Value = "  +10.65%  "
$ws.Range("D51").Value = "8.79"
$ws.Range("E51").Value = "  +9.87%  "

Write-Host "Applied cryptos update"
